$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "66.604.74"
$ws.Cells.Item(2, 5).Value = "  -2.28%  "
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "3.472.44"
$ws.Cells.Item(3, 5).Value = "  -2.35%  "
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "0.999"
$ws.Cells.Item(4, 5).Value = "  -0.16%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "600.67"
$ws.Cells.Item(5, 5).Value = "  -3.04%  "
$ws.Cells.Item(6, 5).Value = "  -4.99%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "3.469.09"
$ws.Cells.Item(7, 5).Value = "  -2.43%  "
$ws.Cells.Item(8, 5).Value = "  +0.14%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.480"
$ws.Cells.Item(9, 5).Value = "  -1.95%  "
$ws.Cells.Item(10, 5).Value = "  -3.06%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "7.55"
$ws.Cells.Item(11, 5).Value = "  +3.31%  "
$ws.Cells.Item(12, 5).Value = "  -3.46%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.0000212"
$ws.Cells.Item(13, 5).Value = "  -3.98%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "4.055.85"
$ws.Cells.Item(14, 5).Value = "  -2.47%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "31.29"
$ws.Cells.Item(15, 5).Value = "  -5.76%  "
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "3.470.59"
$ws.Cells.Item(16, 5).Value = "  -2.26%  "
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "66.680.72"
$ws.Cells.Item(17, 5).Value = "  -2.70%  "
$ws.Cells.Item(18, 5).Value = "  +0.32%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "6.40"
$ws.Cells.Item(19, 5).Value = "  -5.55%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "15.28"
$ws.Cells.Item(20, 5).Value = "  -4.15%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "9.99"
$ws.Cells.Item(21, 5).Value = "  +0.51%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "434.54"
$ws.Cells.Item(22, 5).Value = "  -4.88%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "0.606"
$ws.Cells.Item(23, 5).Value = "  -5.32%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "79.26"
$ws.Cells.Item(24, 5).Value = "  +1.11%  "
$ws.Cells.Item(25, 5).Value = "  -0.07%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "3.609.89"
$ws.Cells.Item(26, 5).Value = "  -2.29%  "
$ws.Cells.Item(27, 5).Value = "  -7.94%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "9.76"
$ws.Cells.Item(28, 5).Value = "  -7.45%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "8.35"
$ws.Cells.Item(29, 5).Value = "  -8.02%  "
$ws.Cells.Item(30, 5).Value = "  -3.41%  "
$ws.Cells.Item(31, 5).Value = "  -6.13%  "
$ws.Cells.Item(32, 5).Value = "  -2.36%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "1.00"
$ws.Cells.Item(33, 5).Value = "  +0.12%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "25.29"
$ws.Cells.Item(34, 5).Value = "  -3.05%  "
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "3.461.90"
$ws.Cells.Item(35, 5).Value = "  -2.51%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "5.94"
$ws.Cells.Item(36, 5).Value = "  -7.09%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "1.79"
$ws.Cells.Item(37, 5).Value = "  -6.45%  "
$ws.Cells.Item(38, 5).Value = "  -0.01%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "7.88"
$ws.Cells.Item(39, 5).Value = "  -4.42%  "
$ws.Cells.Item(40, 5).Value = "  -0.09%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "176.03"
$ws.Cells.Item(41, 5).Value = "  -1.35%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.0879"
$ws.Cells.Item(42, 5).Value = "  -4.24%  "
$ws.Cells.Item(43, 5).Value = "  -11.01%  "
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "5.40"
$ws.Cells.Item(44, 5).Value = "  -3.34%  "
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.891"
$ws.Cells.Item(45, 5).Value = "  -0.50%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "46.33"
$ws.Cells.Item(46, 5).Value = "  -0.50%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "28.77"
$ws.Cells.Item(47, 5).Value = "  -6.86%  "
$ws.Cells.Item(48, 5).Value = "  -8.29%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "7.42"
$ws.Cells.Item(49, 5).Value = "  -4.64%  "
$ws.Cells.Item(50, 5).Value = "  -8.82%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.976"
$ws.Cells.Item(51, 5).Value = "  -4.38%  "
